$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting N:P -> O:Q.
# Excel's default "insert" behaviour copies formatting (incl. width) from
# the column immediately to the left, so mirror that for the column width.
$ws.Columns("N:N").Insert() | Out-Null
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab (was "Edit Repayment
# Schedule" before) and move the selection on it.
$ws.Activate() | Out-Null
$ws.Range("R8").Select() | Out-Null
